# Upload new version with timestamp
# Insert 4 new product rows (alphabetically sorted into the existing list),
# renumber the "م" index column, update the grand total, and bump the
# generated-at timestamp in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the new blank rows.
#    Work on ORIGINAL (pre-edit) row numbers, bottom-to-top, so each
#    insertion doesn't disturb the row numbers used by the next one.
#       - before "نيفيا سوفت كريم 50 مل" (row 42): two new rows
#       - before "فرشه اسنان POWER GOLD كبار" (row 36): one new row
#       - before "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF." (row 19): one new row
# ---------------------------------------------------------------------
$xlShiftDown = -4121

$ws.Rows.Item(42).Insert($xlShiftDown)
$ws.Rows.Item(42).Insert($xlShiftDown)
$ws.Rows.Item(36).Insert($xlShiftDown)
$ws.Rows.Item(19).Insert($xlShiftDown)

# ---------------------------------------------------------------------
# 2) Helper to populate one data row with the same shape as every other
#    product row: A/B merged index, C:G merged name, H:K merged balance,
#    L:M merged reorder flag, N:O merged price, P sell price, Q deal count.
# ---------------------------------------------------------------------
function Fill-ProductRow($row, $index, $name, $balance, $reorder, $price, $sellPrice, $deals, $height) {
    $ws.Rows.Item($row).RowHeight = $height

    $ws.Range("A" + $row + ":B" + $row).Merge() | Out-Null
    $ws.Range("C" + $row + ":G" + $row).Merge() | Out-Null
    $ws.Range("H" + $row + ":K" + $row).Merge() | Out-Null
    $ws.Range("L" + $row + ":M" + $row).Merge() | Out-Null
    $ws.Range("N" + $row + ":O" + $row).Merge() | Out-Null

    $ws.Range("C" + $row + ":G" + $row).NumberFormat = "@"
    $ws.Range("H" + $row + ":K" + $row).NumberFormat = "@"
    $ws.Range("L" + $row + ":M" + $row).NumberFormat = "@"
    $ws.Range("N" + $row + ":O" + $row).NumberFormat = "@"
    $ws.Range("P" + $row).NumberFormat = "@"
    $ws.Range("Q" + $row).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $index
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 8).Value = $balance
    $ws.Cells.Item($row, 12).Value = $reorder
    $ws.Cells.Item($row, 14).Value = $price
    $ws.Cells.Item($row, 16).Value = $sellPrice
    $ws.Cells.Item($row, 17).Value = $deals
}

# ---------------------------------------------------------------------
# 3) Fill in the 4 new rows (final row numbers, after all inserts above).
# ---------------------------------------------------------------------
Fill-ProductRow 19 13 "GOLD PLUS VAG.DOUCHE 250 ML" "0:0" "1" "85.00" "85.0000" "1:0" 25.5
Fill-ProductRow 37 31 "صوفي طويل جدا" "3:0" "0" "50.00" "50.0000" "1:0" 25.5
Fill-ProductRow 44 38 "مناديل بكر فاين" "3:0" "0" "15.00" "15.0000" "1:0" 25.5
Fill-ProductRow 45 39 "مناديل مبلله كبيره" "6:0" "0" "30.00" "30.0000" "1:0" 24.75

# ---------------------------------------------------------------------
# 4) Renumber the "م" column (A) for every data row, 1..40 sequentially.
# ---------------------------------------------------------------------
for ($r = 7; $r -le 46; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 6)
}

# ---------------------------------------------------------------------
# 5) Update the grand total (sum of the sell-price column) and the
#    generated-at timestamp in the footer.
# ---------------------------------------------------------------------
$ws.Cells.Item(47, 16).Value = 1513.435
$ws.Cells.Item(48, 1).Value = "Thursday, 31 July, 2025 3:44 PM"
